$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2597
$ws1.Range("F4").Value = 347
$ws1.Range("F5").Value = 1451
$ws1.Range("F6").Value = 1126
$ws1.Range("F7").Value = 326
$ws1.Range("F13").Value = 8937
$ws1.Range("F19").Value = 611
$ws1.Range("F21").Value = 1167
$ws1.Range("F23").Value = 2062
$ws1.Range("F24").Value = 2138
$ws1.Range("F26").Value = 1842
$ws1.Range("F30").Value = 605
$ws1.Range("F31").Value = 60
$ws1.Range("F32").Value = 134
$ws1.Range("F33").Value = 200
$ws1.Range("F34").Value = 19
$ws1.Range("F38").Value = 458
$ws1.Range("F39").Value = 1320
$ws1.Range("F41").Value = 61
$ws1.Range("F43").Value = 281

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2597
$ws4.Range("F4").Value = 347
$ws4.Range("F5").Value = 1451
$ws4.Range("F7").Value = 1126
$ws4.Range("F8").Value = 326
$ws4.Range("F13").Value = 8937
$ws4.Range("F20").Value = 611
$ws4.Range("F22").Value = 1167
$ws4.Range("F24").Value = 2062
$ws4.Range("F25").Value = 2138
$ws4.Range("F27").Value = 1842
$ws4.Range("F31").Value = 605
$ws4.Range("F32").Value = 60
$ws4.Range("F33").Value = 134
$ws4.Range("F34").Value = 200
$ws4.Range("F35").Value = 19
$ws4.Range("F39").Value = 458
$ws4.Range("F44").Value = 1320
$ws4.Range("F47").Value = 61
$ws4.Range("F49").Value = 281
